# Updated symbol list on Fri Feb 17 07:24:08 UTC 2023 with GitHub Actions
#
# Applies updated Price ("D" column) and Volume(1h) ("E" column) values for
# the cryptos sheet, per the latest scrape. All values are stored as literal
# text (matching the original inlineStr cell type) by prefixing with an
# apostrophe, which forces Excel to keep them as text instead of silently
# re-interpreting numeric- or percent-looking strings as numbers (which would
# corrupt formatting/precision, e.g. turning "310.14" into 310.13999999999999
# or "-3.46%" into a floating point fraction).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.14"
$ws.Range("E2").Value = "'-3.46%"
$ws.Range("D3").Value = "'49.75"
$ws.Range("E3").Value = "'1.30%"
$ws.Range("D4").Value = "'5.166"
$ws.Range("E4").Value = "'-2.38%"
$ws.Range("D5").Value = "'0.07759"
$ws.Range("D6").Value = "'4.520"
$ws.Range("D7").Value = "'1.372"
$ws.Range("E7").Value = "'13.80%"
$ws.Range("D8").Value = "'1.555"
$ws.Range("E8").Value = "'-6.26%"
$ws.Range("D9").Value = "'0.1234"
$ws.Range("E9").Value = "'-6.58%"
$ws.Range("D10").Value = "'0.1975"
$ws.Range("E10").Value = "'0.89%"
$ws.Range("D11").Value = "'0.04722"
$ws.Range("E11").Value = "'6.22%"
$ws.Range("D12").Value = "'0.09464"
$ws.Range("E12").Value = "'-0.85%"
$ws.Range("D13").Value = "'0.1047"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("D14").Value = "'0.001272"
$ws.Range("E14").Value = "'-4.14%"
$ws.Range("D15").Value = "'0.04171"
$ws.Range("E15").Value = "'-3.24%"
$ws.Range("D16").Value = "'0.005795"
$ws.Range("E16").Value = "'-1.41%"
$ws.Range("E17").Value = "'2,014.34%"
$ws.Range("E18").Value = "'-0.87%"
$ws.Range("D19").Value = "'2.241"
$ws.Range("E19").Value = "'-8.04%"
$ws.Range("D20").Value = "'0.3487"
$ws.Range("D21").Value = "'7.905"
$ws.Range("E21").Value = "'-4.40%"
$ws.Range("D22").Value = "'0.1339"
$ws.Range("E22").Value = "'-5.03%"
$ws.Range("E23").Value = "'4.04%"
$ws.Range("E24").Value = "'-2.91%"
$ws.Range("D25").Value = "'0.004025"
$ws.Range("E25").Value = "'-5.32%"
$ws.Range("E26").Value = "'-0.36%"
$ws.Range("D38").Value = "'0.02606"
$ws.Range("E38").Value = "'-4.21%"
$ws.Range("D39").Value = "'0.05862"
$ws.Range("E39").Value = "'4.75%"
$ws.Range("E40").Value = "'69.65%"
$ws.Range("D41").Value = "'0.007905"
$ws.Range("E41").Value = "'2.62%"
$ws.Range("E42").Value = "'-1.59%"
$ws.Range("D43").Value = "'0.008445"
$ws.Range("E43").Value = "'9.50%"
$ws.Range("D44").Value = "'0.007659"
$ws.Range("E44").Value = "'-5.46%"
$ws.Range("D45").Value = "'0.3402"
$ws.Range("E45").Value = "'6.51%"
$ws.Range("D46").Value = "'0.00007031"
$ws.Range("E46").Value = "'0.27%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.35%"
$ws.Range("D48").Value = "'0.05129"
$ws.Range("E48").Value = "'-16.36%"
$ws.Range("D49").Value = "'0.002622"
$ws.Range("E49").Value = "'-34.66%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.35%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.35%"
